$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '37.357.18'
$ws.Range('E2').Value = '  +2.10%  '

# Row 3
$ws.Range('D3').Value = '2.089.23'
$ws.Range('E3').Value = '  +0.06%  '

# Row 4
$ws.Range('E4').Value = '  -0.01%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '252.04'
$ws.Range('E5').Value = '  +1.84%  '

# Row 6
$ws.Range('E6').Value = '  -0.05%  '

# Row 8
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '54.50'
$ws.Range('E8').Value = '  +20.82%  '

# Row 9
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '62.45'
$ws.Range('E9').Value = '  +3.02%  '

# Row 10
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.382'
$ws.Range('E10').Value = '  +4.88%  '

# Row 11
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0755'
$ws.Range('E11').Value = '  +4.33%  '

# Row 12
$ws.Range('E12').Value = '  +7.41%  '

# Row 13
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '15.52'
$ws.Range('E13').Value = '  +6.84%  '

# Row 14
$ws.Range('B14').Value = 'Polygon'
$ws.Range('C14').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.870'
$ws.Range('E14').Value = '  +5.31%  '

# Row 15
$ws.Range('B15').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C15').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D15').Value = '2.393.56'
$ws.Range('E15').Value = '  +0.32%  '

# Row 16
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '5.27'
$ws.Range('E16').Value = '  +6.98%  '

# Row 17
$ws.Range('D17').Value = '2.093.50'
$ws.Range('E17').Value = '  +0.78%  '

# Row 18
$ws.Range('D18').Value = '37.266.82'
$ws.Range('E18').Value = '  +1.79%  '

# Row 19
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '73.34'
$ws.Range('E19').Value = '  +2.15%  '

# Row 20
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '14.50'
$ws.Range('E20').Value = '  +13.84%  '

# Row 21
$ws.Range('E21').Value = '  +4.89%  '

# Row 22
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '241.30'
$ws.Range('E22').Value = '  +1.23%  '

# Row 23
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.27'
$ws.Range('E23').Value = '  +6.78%  '

# Row 24
$ws.Range('E24').Value = '  -0.04%  '

# Row 25
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.47'
$ws.Range('E25').Value = '  +0.21%  '

# Row 26
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '172.30'
$ws.Range('E26').Value = '  +1.76%  '

# Row 27
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.29'
$ws.Range('E27').Value = '  +4.88%  '

# Row 28
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '20.92'
$ws.Range('E28').Value = '  +2.09%  '

# Row 29
$ws.Range('E29').Value = '  +3.91%  '

# Row 30
$ws.Range('E30').Value = '  +2.24%  '

# Row 31
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '23.39'
$ws.Range('E31').Value = '  +4.69%  '

# Row 32
$ws.Range('E32').Value = '  +21.88%  '

# Row 33
$ws.Range('E33').Value = '  +4.02%  '

# Row 34
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.0628'

# Row 35
$ws.Range('B35').Value = 'InternetComputer(DFINITY)'
$ws.Range('C35').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '4.33'
$ws.Range('E35').Value = '  +8.21%  '

# Row 36
$ws.Range('B36').Value = 'Kaspa'
$ws.Range('C36').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.0904'
$ws.Range('E36').Value = '  -0.15%  '

# Row 37
$ws.Range('E37').Value = '  +0.05%  '

# Row 38
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.27'
$ws.Range('E38').Value = '  -1.18%  '

# Row 39
$ws.Range('E39').Value = '  -3.16%  '

# Row 40
$ws.Range('E40').Value = '  +1.19%  '

# Row 41
$ws.Range('E41').Value = '  +6.04%  '

# Row 42
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '17.92'
$ws.Range('E42').Value = '  +12.87%  '

# Row 43
$ws.Range('E43').Value = '  +2.77%  '

# Row 44
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '99.85'
$ws.Range('E44').Value = '  +2.11%  '

# Row 45
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0971'
$ws.Range('E45').Value = '  +18.09%  '

# Row 46
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '4.41'
$ws.Range('E46').Value = '  +122.16%  '

# Row 47
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.80'
$ws.Range('E47').Value = '  +0.56%  '

# Row 48
$ws.Range('D48').Value = '1.333.06'
$ws.Range('E48').Value = '  -0.24%  '

# Row 49
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.95'
$ws.Range('E49').Value = '  +4.10%  '

# Row 50
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.36'
$ws.Range('E50').Value = '  +6.52%  '

# Row 51
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '7.01'
$ws.Range('E51').Value = '  +12.64%  '
